# Auto-generated Excel COM-interop script
# Applies updated crypto price/volume data, and swaps Maker/ApeXProtocol rows 47-48
# (matches the commit 'Updated cryptos list ... with GitHub Actions')

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '50.873.15'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -2.41%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.902.69'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -2.51%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '371.91'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +5.06%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '101.58'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -5.39%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.541'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -3.90%  '
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -5.02%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.74'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -3.96%  '
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +0.36%  '
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -2.96%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '18.20'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -5.51%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.360.40'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -2.49%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.33'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -3.80%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.902.71'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -2.29%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.919'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -8.05%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '50.812.45'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -2.63%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.23'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -7.12%  '
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -4.27%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.83'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -5.70%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.0₃0938'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -3.72%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '68.09'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -2.13%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '258.94'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -1.94%  '
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -1.96%  '
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -5.40%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '4.07'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -4.98%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.03'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -7.34%  '
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -6.53%  '
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +1.38%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '9.82'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -4.60%  '
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -2.75%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '51.14'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +0.77%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '33.93'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -6.94%  '
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +0.28%  '
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -4.99%  '
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -7.09%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '16.92'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -5.46%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.57'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -5.69%  '
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -7.17%  '
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -4.29%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '118.91'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -2.13%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '21.76'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -4.20%  '
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -1.70%  '
$ws.Range('B47').NumberFormat = "@"
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').NumberFormat = "@"
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.31'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -1.66%  '
$ws.Range('B48').NumberFormat = "@"
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').NumberFormat = "@"
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.006.29'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -5.26%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '3.13'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -6.98%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '3.194.25'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -2.33%  '
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -2.47%  '

Write-Host "Applied 83 cell updates."
